$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.888.23"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "3.165.67"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'616.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.52%  "
$ws.Range("D6").Value = "'148.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.72%  "
$ws.Range("D8").Value = "3.162.19"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("D11").Value = "'5.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.95%  "
$ws.Range("D12").Value = "'0.474"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("D13").Value = "'0.0000260"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").Value = "'35.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.39%  "
$ws.Range("D15").Value = "3.677.79"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("E16").Value = "  +2.68%  "
$ws.Range("D17").Value = "64.836.97"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "3.162.72"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").Value = "'6.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.49%  "
$ws.Range("D20").Value = "'482.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").Value = "'14.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("D22").Value = "'0.722"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("D23").Value = "'8.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.33%  "
$ws.Range("D24").Value = "'13.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.45%  "
$ws.Range("D25").Value = "'84.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  -2.78%  "
$ws.Range("E28").Value = "  -2.07%  "
$ws.Range("D29").Value = "'6.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.06%  "
$ws.Range("D30").Value = "'0.119"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.55%  "
$ws.Range("E31").Value = "  -7.57%  "
$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").Value = "'2.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").Value = "'26.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.27%  "
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("D36").Value = "0.0₃0786"
$ws.Range("E36").Value = "  +5.20%  "
$ws.Range("D37").Value = "'6.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.02%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").Value = "'53.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.37%  "
$ws.Range("D40").Value = "'461.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("E41").Value = "  -0.39%  "
$ws.Range("E42").Value = "  -4.69%  "
$ws.Range("D43").Value = "'8.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.53%  "
$ws.Range("D44").Value = "2.848.77"
$ws.Range("E44").Value = "  -1.64%  "
$ws.Range("E45").Value = "  -4.37%  "
$ws.Range("D46").Value = "'0.270"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("E47").Value = "  +4.99%  "
$ws.Range("D48").Value = "'26.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.62%  "
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("E50").Value = "  -1.34%  "
$ws.Range("D51").Value = "'120.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.81%  "
